$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the title banner in A1 to the new version/date.
$ws.Range("A1").Value = "sp_AskBrent Check ID List - v20 2016-01-01"

# 2. Append the new check row (27) at the bottom of the table.
$ws.Cells.Item(31, 1).Value = 27
$ws.Cells.Item(31, 2).Value = 1
$ws.Cells.Item(31, 3).Value = "Outdated sp_AskBrent"
$ws.Cells.Item(31, 4).Value = "sp_AskBrent is Over 6 Months Old"
$ws.Cells.Item(31, 5).Value = "http://BrentOzar.com/askbrent/"

# Hyperlink the URL cell (Excel's own link-target auto-correct lower-cases
# the host, matching every other brentozar.com link already in this sheet),
# then restore the same hyperlink look-and-feel used by the rest of column E
# (re-using the existing "Hyperlink" style rather than leaving the freshly
# minted one that Hyperlinks.Add applies automatically).
$ws.Hyperlinks.Add($ws.Cells.Item(31, 5), "http://brentozar.com/askbrent/")
$ws.Cells.Item(31, 5).Style = $ws.Cells.Item(30, 5).Style

# 3. Move the active selection in the frozen bottom-right pane to B5.
$ws.Range("B5").Select()

Write-Host "done"
